# Apply the target changes to ProjectIndex.xlsx
$wb = $excel.ActiveWorkbook

# Sheet "ABC Notes" (2nd tab) was the active/selected sheet; it should no
# longer be the tab-selected sheet.
$wsNotes = $wb.Worksheets.Item("ABC Notes")

# Sheet "Excel Notes" (3rd tab) becomes the active/selected sheet,
# with selection/active cell moved from B3 to C2.
$wsExcelNotes = $wb.Worksheets.Item("Excel Notes")

# Clear the value in C2 on the "Excel Notes" sheet (removes the <c> element).
$wsExcelNotes.Range("C2").ClearContents()

# Update the selected cell on "Excel Notes" to C2.
$wsExcelNotes.Range("C2").Select()

# Activate "Excel Notes" so it becomes the active tab / tabSelected sheet,
# and the workbook's active tab pointer moves to this sheet.
$wsExcelNotes.Activate()
